# Auto-generated edit script
# Updates numeric price/profit columns (H-N) across multiple sheets
# per scheduled data-refresh run.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1209862.2
$ws.Range("J17").Value = 1245973.2
$ws.Range("L17").Value = 3737919.6
$ws.Range("N17").Value = -3738255.6
$ws.Range("H64").Value = 3456.0588
$ws.Range("I64").Value = 3418.75
$ws.Range("K64").Value = 3418.75
$ws.Range("M64").Value = -3170.75
$ws.Range("H67").Value = 3456.0588
$ws.Range("I67").Value = 3418.75
$ws.Range("K67").Value = 3418.75
$ws.Range("M67").Value = -2560.75
$ws.Range("H116").Value = 2145
$ws.Range("I116").Value = 1818.5714
$ws.Range("J116").Value = 2471.4285
$ws.Range("K116").Value = 1818.5714
$ws.Range("L116").Value = 2471.4285
$ws.Range("M116").Value = 1623.4286
$ws.Range("N116").Value = -9355.4285
$ws.Range("H129").Value = 3180.9524
$ws.Range("J129").Value = 3180.9524
$ws.Range("L129").Value = 9542.8572
$ws.Range("N129").Value = -19542.8572
$ws.Range("H132").Value = 2436.07
$ws.Range("I132").Value = 2254.875
$ws.Range("J132").Value = 3402.4443
$ws.Range("K132").Value = 6764.625
$ws.Range("L132").Value = 10207.3329
$ws.Range("M132").Value = -4234.625
$ws.Range("N132").Value = -15267.3329
$ws.Range("H137").Value = 5715241.5
$ws.Range("I137").Value = 943.5161000000001
$ws.Range("J137").Value = 50001050
$ws.Range("K137").Value = 2830.5483
$ws.Range("L137").Value = 150003150
$ws.Range("M137").Value = -280.5483000000004
$ws.Range("N137").Value = -150008250

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H19").Value = 3877.25
$ws.Range("I19").Value = 5750
$ws.Range("J19").Value = 2004.5
$ws.Range("K19").Value = 5750
$ws.Range("L19").Value = 2004.5
$ws.Range("M19").Value = -5521
$ws.Range("N19").Value = -2462.5
$ws.Range("H45").Value = 1482.8572
$ws.Range("I45").Value = 1443.2354
$ws.Range("J45").Value = 1651.25
$ws.Range("K45").Value = 1443.2354
$ws.Range("L45").Value = 1651.25
$ws.Range("M45").Value = -1066.2354
$ws.Range("N45").Value = -2405.25
$ws.Range("H74").Value = 11630570
$ws.Range("I74").Value = 19232568
$ws.Range("J74").Value = 3984.5881
$ws.Range("K74").Value = 19232568
$ws.Range("L74").Value = 3984.5881
$ws.Range("M74").Value = -19231694
$ws.Range("N74").Value = -5732.5881
$ws.Range("H77").Value = 11630570
$ws.Range("I77").Value = 19232568
$ws.Range("J77").Value = 3984.5881
$ws.Range("K77").Value = 96162840
$ws.Range("L77").Value = 19922.9405
$ws.Range("M77").Value = -96158472
$ws.Range("N77").Value = -28658.9405
$ws.Range("H102").Value = 1876
$ws.Range("I102").Value = 1702
$ws.Range("J102").Value = 2050
$ws.Range("K102").Value = 1702
$ws.Range("L102").Value = 2050
$ws.Range("M102").Value = -80
$ws.Range("N102").Value = -5294
$ws.Range("H132").Value = 8623419
$ws.Range("I132").Value = 15627382
$ws.Range("K132").Value = 46882146
$ws.Range("M132").Value = -46879616

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3837.0356
$ws.Range("I134").Value = 2539.8572
$ws.Range("J134").Value = 7728.5713
$ws.Range("K134").Value = 7619.571599999999
$ws.Range("L134").Value = 23185.7139
$ws.Range("M134").Value = -5084.571599999999
$ws.Range("N134").Value = -28255.7139

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H26").Value = 13882.875
$ws.Range("I26").Value = 2200
$ws.Range("J26").Value = 33354.332
$ws.Range("K26").Value = 2200
$ws.Range("L26").Value = 33354.332
$ws.Range("M26").Value = -1913
$ws.Range("N26").Value = -33928.332
$ws.Range("H31").Value = 6669932
$ws.Range("I31").Value = 3338.75
$ws.Range("J31").Value = 166668180
$ws.Range("K31").Value = 3338.75
$ws.Range("L31").Value = 166668180
$ws.Range("M31").Value = -3043.75
$ws.Range("N31").Value = -166668770
$ws.Range("H34").Value = 6669932
$ws.Range("I34").Value = 3338.75
$ws.Range("J34").Value = 166668180
$ws.Range("K34").Value = 3338.75
$ws.Range("L34").Value = 166668180
$ws.Range("M34").Value = -3136.75
$ws.Range("N34").Value = -166668584
$ws.Range("H56").Value = 0
$ws.Range("J56").Value = 0
$ws.Range("L56").Value = 0
$ws.Range("N56").ClearContents()
$ws.Range("H58").Value = 1459.2128
$ws.Range("I58").Value = 752.1111
$ws.Range("J58").Value = 2413.8
$ws.Range("K58").Value = 752.1111
$ws.Range("L58").Value = 2413.8
$ws.Range("M58").Value = -549.1111
$ws.Range("N58").Value = -2819.8
$ws.Range("H62").Value = 3111
$ws.Range("J62").Value = 4214.2
$ws.Range("L62").Value = 4214.2
$ws.Range("N62").Value = -5462.2
$ws.Range("H65").Value = 3111
$ws.Range("J65").Value = 4214.2
$ws.Range("L65").Value = 21071
$ws.Range("N65").Value = -27311
$ws.Range("H86").Value = 2545.9375
$ws.Range("I86").Value = 2552.5
$ws.Range("K86").Value = 2552.5
$ws.Range("M86").Value = -1429.5
$ws.Range("H89").Value = 2545.9375
$ws.Range("I89").Value = 2552.5
$ws.Range("K89").Value = 12762.5
$ws.Range("M89").Value = -7146.5
$ws.Range("H132").Value = 11906972
$ws.Range("I132").Value = 15153189
$ws.Range("J132").Value = 4178.3335
$ws.Range("K132").Value = 45459567
$ws.Range("L132").Value = 12535.0005
$ws.Range("M132").Value = -45457037
$ws.Range("N132").Value = -17595.0005
$ws.Range("H134").Value = 1318.4849
$ws.Range("I134").Value = 1431.36
$ws.Range("J134").Value = 965.75
$ws.Range("K134").Value = 4294.08
$ws.Range("L134").Value = 2897.25
$ws.Range("M134").Value = -1759.08
$ws.Range("N134").Value = -7967.25
$ws.Range("H136").Value = 1459.2128
$ws.Range("I136").Value = 752.1111
$ws.Range("J136").Value = 2413.8
$ws.Range("K136").Value = 2256.3333
$ws.Range("L136").Value = 7241.400000000001
$ws.Range("M136").Value = 293.6667000000002
$ws.Range("N136").Value = -12341.4
$ws.Range("H140").Value = 38096
$ws.Range("J140").Value = 38096
$ws.Range("L140").Value = 38096
$ws.Range("N140").Value = -48456

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 142.46666
$ws.Range("I23").Value = 94.2
$ws.Range("K23").Value = 282.6
$ws.Range("M23").Value = -47.60000000000002
$ws.Range("H33").Value = 42857228
$ws.Range("I33").Value = 36363724
$ws.Range("K33").Value = 218182344
$ws.Range("M33").Value = -218182061
$ws.Range("H97").Value = 21641.2
$ws.Range("I97").Value = 26051.5
$ws.Range("J97").Value = 4000
$ws.Range("K97").Value = 78154.5
$ws.Range("L97").Value = 12000
$ws.Range("M97").Value = -77658.5
$ws.Range("N97").Value = -12992
$ws.Range("H123").Value = 6166.625
$ws.Range("I123").Value = 2000
$ws.Range("J123").Value = 7555.5
$ws.Range("K123").Value = 6000
$ws.Range("L123").Value = 22666.5
$ws.Range("M123").Value = -3550
$ws.Range("N123").Value = -27566.5
$ws.Range("H131").Value = 833.08246
$ws.Range("J131").Value = 858.9888999999999
$ws.Range("L131").Value = 2576.9667
$ws.Range("N131").Value = -12656.9667
$ws.Range("H132").Value = 867.9167
$ws.Range("I132").Value = 685.7143
$ws.Range("J132").Value = 1123
$ws.Range("K132").Value = 6171.428699999999
$ws.Range("L132").Value = 10107
$ws.Range("M132").Value = -3641.428699999999
$ws.Range("N132").Value = -15167
$ws.Range("H137").Value = 3494.0789
$ws.Range("I137").Value = 1571.5385
$ws.Range("J137").Value = 4493.8
$ws.Range("K137").Value = 4714.6155
$ws.Range("L137").Value = 13481.4
$ws.Range("M137").Value = 385.3845000000001
$ws.Range("N137").Value = -23681.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3456.0833
$ws.Range("I102").Value = 3900.3157
$ws.Range("K102").Value = 3900.3157
$ws.Range("M102").Value = -2278.3157
$ws.Range("H113").Value = 41105.52
$ws.Range("I113").Value = 84351.914
$ws.Range("J113").Value = 1185.7693
$ws.Range("K113").Value = 84351.914
$ws.Range("L113").Value = 1185.7693
$ws.Range("M113").Value = -82181.914
$ws.Range("N113").Value = -5525.7693
$ws.Range("H122").Value = 4169146.5
$ws.Range("J122").Value = 3895
$ws.Range("L122").Value = 11685
$ws.Range("N122").Value = -16585

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 1682.85
$ws.Range("I68").Value = 1882.8334
$ws.Range("K68").Value = 1882.8334
$ws.Range("M68").Value = -1133.8334
$ws.Range("H71").Value = 1682.85
$ws.Range("I71").Value = 1882.8334
$ws.Range("K71").Value = 9414.166999999999
$ws.Range("M71").Value = -5670.166999999999
$ws.Range("H136").Value = 12199945
$ws.Range("I136").Value = 20002038
$ws.Range("J136").Value = 9175.625
$ws.Range("K136").Value = 60006114
$ws.Range("L136").Value = 27526.875
$ws.Range("M136").Value = -60003564
$ws.Range("N136").Value = -32626.875

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1290.0333
$ws.Range("I100").Value = 1735.5714
$ws.Range("J100").Value = 900.1875
$ws.Range("K100").Value = 3471.1428
$ws.Range("L100").Value = 1800.375
$ws.Range("M100").Value = -2930.1428
$ws.Range("N100").Value = -2882.375
$ws.Range("H126").Value = 2768.5417
$ws.Range("I126").Value = 1259.0476
$ws.Range("K126").Value = 3777.142800000001
$ws.Range("M126").Value = -1307.142800000001
$ws.Range("H132").Value = 2619.2942
$ws.Range("I132").Value = 1488.6364
$ws.Range("J132").Value = 4692.1665
$ws.Range("K132").Value = 4465.9092
$ws.Range("L132").Value = 14076.4995
$ws.Range("M132").Value = -1935.9092
$ws.Range("N132").Value = -19136.4995
$ws.Range("H136").Value = 1149.7368
$ws.Range("I136").Value = 980.9375
$ws.Range("J136").Value = 2050
$ws.Range("K136").Value = 2942.8125
$ws.Range("L136").Value = 6150
$ws.Range("M136").Value = -392.8125
$ws.Range("N136").Value = -11250
